$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "44.036.57"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "2.244.06"
$ws.Range("E3").Value = "  +2.07%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "99.43"
$ws.Range("E5").Value = "  +19.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "271.26"
$ws.Range("E6").Value = "  +4.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  +1.18%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  +7.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.33"
$ws.Range("E10").Value = "  +7.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0947"
$ws.Range("E11").Value = "  +3.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.36"
$ws.Range("E12").Value = "  +17.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.38"
$ws.Range("E14").Value = "  +7.09%  "
$ws.Range("D15").Value = "2.541.25"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.819"
$ws.Range("E16").Value = "  +4.54%  "
$ws.Range("D17").Value = "2.246.44"
$ws.Range("E17").Value = "  +2.56%  "
$ws.Range("D18").Value = "44.013.70"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("E19").Value = "  +2.53%  "
$ws.Range("E20").Value = "  +5.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.26"
$ws.Range("E21").Value = "  +2.35%  "
$ws.Range("E22").Value = "  -2.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.08"
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.24"
$ws.Range("E24").Value = "  +2.67%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.63"
$ws.Range("E26").Value = "  +9.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.52"
$ws.Range("E27").Value = "  +13.25%  "
$ws.Range("E28").Value = "  +3.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.76"
$ws.Range("E29").Value = "  +2.98%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.88"
$ws.Range("E31").Value = "  -0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0934"
$ws.Range("E32").Value = "  +8.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.04"
$ws.Range("E33").Value = "  +3.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.61"
$ws.Range("E34").Value = "  +4.96%  "
$ws.Range("E35").Value = "  +1.69%  "
$ws.Range("E36").Value = "  +2.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0353"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("E38").Value = "  -3.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.74"
$ws.Range("E39").Value = "  +33.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.249"
$ws.Range("E40").Value = "  +24.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.93"
$ws.Range("E41").Value = "  +3.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.19"
$ws.Range("E42").Value = "  +4.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.85"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.45"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("E45").Value = "  +4.06%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.50"
$ws.Range("E46").Value = "  +2.60%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.79"
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("E48").Value = "  +4.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.19"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.439"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("E51").Value = "  +1.14%  "
